# Daily attendance processing - 2026-01-26 00:00:41
# Normalizes the "Recorded By" column (G) so that any value whose first
# comma-separated token is "System" has its first and last comma-separated
# tokens swapped (e.g. "System, x@y.com" -> "x@y.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Length -gt 0) {
        $parts = $val -split ", "
        if ($parts.Length -ge 2 -and $parts[0] -eq "System") {
            $first = $parts[0]
            $last = $parts[$parts.Length - 1]
            $parts[0] = $last
            $parts[$parts.Length - 1] = $first
            $cell.Value = [string]::Join(", ", $parts)
        }
    }
}
